$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Individualised Learner Record rows with the latest data release
# details (row 5: "Further education and skills achievements", row 6:
# "...by sector subject area").
$ws.Range("B5").Value = "<a href='https://explore-education-statistics.service.gov.uk/data-catalogue/further-education-and-skills/2021-22'>Individualised Learner Record</a>"
$ws.Range("C5").Value = "Aug 2021 – Jul 2022 (24/11/22)"
$ws.Range("D5").Value = "Aug 2022 – Jan 2023 (Mar 23)"

$ws.Range("B6").Value = "<a href='https://explore-education-statistics.service.gov.uk/data-tables/permalink/cae2bcbb-e385-4da7-8d7b-08dacbbccc68'>Individualised Learner Record</a>"
$ws.Range("C6").Value = "Aug 2021 – Jul 2022 (24/11/22)"
$ws.Range("D6").Value = "Aug 2022 – Jan 2023 (Mar 23)"

# Update the window/selection state to match the latest save.
$ws.Range("B6").Select()
